$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: "FRI Jan 26" and " 14:07:37 PST 2018" were two separate
# runs; the edit merges them into a single run's text. Functionally
# this is just a text replace across that span.
# ---------------------------------------------------------------------
$find = $d.Content.Find
$find.ClearFormatting()
[void]$find.Execute("FRI Jan 26 14:07:37 PST 2018", $true, $false, $false, $false, $false, $true, 1, $false, "FRI Jan 26 14:07:37 PST 2018", 2)

# ---------------------------------------------------------------------
# Change 2: append a new purchase-record block (03/02/2018 HARISH CHICK
# IN) right after the paragraph ending "Amount balance ... - 129002.0",
# before the pre-existing trailing blank paragraphs.
# ---------------------------------------------------------------------

# Locate the paragraph that holds "- 129002.0".
$paras = $d.Paragraphs
$targetIndex = -1
for ($i = 1; $i -le $paras.Count; $i++) {
    if ($paras.Item($i).Range.Text -like "*129002.0*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq -1) {
    throw "Could not find the '- 129002.0' paragraph"
}

# Create a fresh blank paragraph right after it, then replace that
# blank paragraph's content with the full new block via InsertXML so
# the run/tab structure matches exactly (tabs as <w:tab/> elements,
# bold runs, PlainText style, Courier New font).
$target = $paras.Item($targetIndex)
$endRng = $target.Range
$endRng.Collapse(0)
[void]$endRng.InsertParagraphAfter()

$newParaRng = $d.Paragraphs.Item($targetIndex + 1).Range

$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$rFonts = '<w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>'

$xml = "<w:p $ns><w:pPr><w:pStyle w:val=`"PlainText`"/><w:rPr>$rFonts<w:b/></w:rPr></w:pPr></w:p>" +
  "<w:p $ns><w:pPr><w:pStyle w:val=`"PlainText`"/><w:rPr>$rFonts</w:rPr></w:pPr>" +
  "<w:r><w:rPr>$rFonts</w:rPr><w:t>FRI Feb 02</w:t></w:r>" +
  "<w:r><w:rPr>$rFonts</w:rPr><w:t xml:space=`"preserve`"> 13:26:18 PST 2018</w:t></w:r></w:p>" +
  "<w:p $ns><w:pPr><w:pStyle w:val=`"PlainText`"/><w:rPr>$rFonts</w:rPr></w:pPr>" +
  "<w:r><w:rPr>$rFonts</w:rPr><w:t>Person Name</w:t></w:r>" +
  "<w:r><w:rPr>$rFonts</w:rPr><w:tab/></w:r>" +
  "<w:r><w:rPr>$rFonts</w:rPr><w:tab/></w:r>" +
  "<w:r><w:rPr>$rFonts</w:rPr><w:tab/></w:r>" +
  "<w:r><w:rPr>$rFonts</w:rPr><w:tab/><w:t>- NNA</w:t></w:r></w:p>" +
  "<w:p $ns><w:pPr><w:pStyle w:val=`"PlainText`"/><w:rPr>$rFonts</w:rPr></w:pPr>" +
  "<w:r><w:rPr>$rFonts</w:rPr><w:t>Bill number</w:t></w:r>" +
  "<w:r><w:rPr>$rFonts</w:rPr><w:tab/></w:r>" +
  "<w:r><w:rPr>$rFonts</w:rPr><w:tab/></w:r>" +
  "<w:r><w:rPr>$rFonts</w:rPr><w:tab/></w:r>" +
  "<w:r><w:rPr>$rFonts</w:rPr><w:tab/><w:t>- 1178</w:t></w:r></w:p>" +
  "<w:p $ns><w:pPr><w:pStyle w:val=`"PlainText`"/><w:rPr>$rFonts</w:rPr></w:pPr>" +
  "<w:r><w:rPr>$rFonts</w:rPr><w:t>---------------------------------------------------------------</w:t></w:r></w:p>" +
  "<w:p $ns><w:pPr><w:pStyle w:val=`"PlainText`"/><w:rPr>$rFonts</w:rPr></w:pPr>" +
  "<w:r><w:rPr>$rFonts</w:rPr><w:t>Item Name</w:t></w:r>" +
  "<w:r><w:rPr>$rFonts</w:rPr><w:tab/></w:r>" +
  "<w:r><w:rPr>$rFonts</w:rPr><w:tab/></w:r>" +
  "<w:r><w:rPr>$rFonts</w:rPr><w:tab/></w:r>" +
  "<w:r><w:rPr>$rFonts</w:rPr><w:tab/><w:t>- CARROT</w:t></w:r></w:p>" +
  "<w:p $ns><w:pPr><w:pStyle w:val=`"PlainText`"/><w:rPr>$rFonts</w:rPr></w:pPr>" +
  "<w:r><w:rPr>$rFonts</w:rPr><w:t>Number of Pockets</w:t></w:r>" +
  "<w:r><w:rPr>$rFonts</w:rPr><w:tab/></w:r>" +
  "<w:r><w:rPr>$rFonts</w:rPr><w:tab/></w:r>" +
  "<w:r><w:rPr>$rFonts</w:rPr><w:tab/><w:t>- 6</w:t></w:r></w:p>" +
  "<w:p $ns><w:pPr><w:pStyle w:val=`"PlainText`"/><w:rPr>$rFonts</w:rPr></w:pPr>" +
  "<w:r><w:rPr>$rFonts</w:rPr><w:t>Number of KGs</w:t></w:r>" +
  "<w:r><w:rPr>$rFonts</w:rPr><w:tab/></w:r>" +
  "<w:r><w:rPr>$rFonts</w:rPr><w:tab/></w:r>" +
  "<w:r><w:rPr>$rFonts</w:rPr><w:tab/><w:t>- 590</w:t></w:r></w:p>" +
  "<w:p $ns><w:pPr><w:pStyle w:val=`"PlainText`"/><w:rPr>$rFonts</w:rPr></w:pPr>" +
  "<w:r><w:rPr>$rFonts</w:rPr><w:t>Rate</w:t></w:r>" +
  "<w:r><w:rPr>$rFonts</w:rPr><w:tab/></w:r>" +
  "<w:r><w:rPr>$rFonts</w:rPr><w:tab/></w:r>" +
  "<w:r><w:rPr>$rFonts</w:rPr><w:tab/></w:r>" +
  "<w:r><w:rPr>$rFonts</w:rPr><w:tab/></w:r>" +
  "<w:r><w:rPr>$rFonts</w:rPr><w:tab/><w:t>- 10</w:t></w:r></w:p>" +
  "<w:p $ns><w:pPr><w:pStyle w:val=`"PlainText`"/><w:rPr>$rFonts</w:rPr></w:pPr>" +
  "<w:r><w:rPr>$rFonts</w:rPr><w:t>Transport &amp; Miscellaneous</w:t></w:r>" +
  "<w:r><w:rPr>$rFonts</w:rPr><w:tab/><w:t>- 60</w:t></w:r></w:p>" +
  "<w:p $ns><w:pPr><w:pStyle w:val=`"PlainText`"/><w:rPr>$rFonts</w:rPr></w:pPr>" +
  "<w:r><w:rPr>$rFonts</w:rPr><w:t>Total Price</w:t></w:r>" +
  "<w:r><w:rPr>$rFonts</w:rPr><w:tab/></w:r>" +
  "<w:r><w:rPr>$rFonts</w:rPr><w:tab/></w:r>" +
  "<w:r><w:rPr>$rFonts</w:rPr><w:tab/></w:r>" +
  "<w:r><w:rPr>$rFonts</w:rPr><w:tab/><w:t>- 5960.0</w:t></w:r></w:p>" +
  "<w:p $ns><w:pPr><w:pStyle w:val=`"PlainText`"/><w:rPr>$rFonts<w:b/></w:rPr></w:pPr>" +
  "<w:r><w:rPr>$rFonts<w:b/></w:rPr><w:t>Amount balance</w:t></w:r>" +
  "<w:r><w:rPr>$rFonts<w:b/></w:rPr><w:tab/></w:r>" +
  "<w:r><w:rPr>$rFonts<w:b/></w:rPr><w:tab/></w:r>" +
  "<w:r><w:rPr>$rFonts<w:b/></w:rPr><w:tab/><w:t>- 134962.0</w:t></w:r></w:p>" +
  "<w:p $ns><w:pPr><w:pStyle w:val=`"PlainText`"/><w:rPr>$rFonts<w:b/></w:rPr></w:pPr></w:p>" +
  "<w:p $ns><w:pPr><w:pStyle w:val=`"PlainText`"/><w:rPr>$rFonts<w:b/></w:rPr></w:pPr></w:p>"

[void]$newParaRng.InsertXML($xml)
